$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update period headers (row 8) and publish-date headers (row 9) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-04-19 (8)"
$ws.Range("E9").Value = "1400-04-14 (8)"
$ws.Range("F9").Value = "1401-04-15 (10)"
$ws.Range("G9").Value = "1402-02-27 (9)"
$ws.Range("H9").Value = "1402-02-27 (2)"

# --- Shift financial data one column to the left (drop oldest period,
#     add newest period 1401/12 in column H) for every data row ---
$ws.Range("D12").Value = 6391857
$ws.Range("E12").Value = 6288642
$ws.Range("F12").Value = 10645362
$ws.Range("G12").Value = 9387301
$ws.Range("H12").Value = 12591159

$ws.Range("D13").Value = 5586658
$ws.Range("E13").Value = 9401658
$ws.Range("F13").Value = 3633000
$ws.Range("G13").Value = 10029498
$ws.Range("H13").Value = 1153779

$ws.Range("D14").Value = 15221995
$ws.Range("E14").Value = 10957104
$ws.Range("F14").Value = 25664499
$ws.Range("G14").Value = 26963306
$ws.Range("H14").Value = 32851840

$ws.Range("D15").Value = 3919066
$ws.Range("E15").Value = 4774835
$ws.Range("F15").Value = 7860309
$ws.Range("G15").Value = 14176468
$ws.Range("H15").Value = 14419573

$ws.Range("D16").Value = 1187981
$ws.Range("E16").Value = 1017567
$ws.Range("F16").Value = 2171160
$ws.Range("G16").Value = 808725
$ws.Range("H16").Value = 941109

$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 0
$ws.Range("G17").Value = 0
$ws.Range("H17").Value = 0

$ws.Range("D18").Value = 32307557
$ws.Range("E18").Value = 32439806
$ws.Range("F18").Value = 49974330
$ws.Range("G18").Value = 61365298
$ws.Range("H18").Value = 61957460

$ws.Range("D19").Value = 0
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

$ws.Range("D20").Value = 184836
$ws.Range("E20").Value = 184836
$ws.Range("F20").Value = 220191
$ws.Range("G20").Value = 220191
$ws.Range("H20").Value = 220191

$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 0
$ws.Range("H21").Value = 0

$ws.Range("D22").Value = 6890470
$ws.Range("E22").Value = 8091192
$ws.Range("F22").Value = 9665642
$ws.Range("G22").Value = 13160553
$ws.Range("H22").Value = 15994145

$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 11100

$ws.Range("D24").Value = "-"
$ws.Range("E24").Value = "-"
$ws.Range("F24").Value = "-"
$ws.Range("G24").Value = "-"
$ws.Range("H24").Value = "-"

$ws.Range("D25").Value = 4810423
$ws.Range("E25").Value = 3524957
$ws.Range("F25").Value = 540525
$ws.Range("G25").Value = 365505
$ws.Range("H25").Value = 541612

$ws.Range("D26").Value = 11885729
$ws.Range("E26").Value = 11800985
$ws.Range("F26").Value = 10426358
$ws.Range("G26").Value = 13746249
$ws.Range("H26").Value = 16767048

$ws.Range("D27").Value = 44193286
$ws.Range("E27").Value = 44240791
$ws.Range("F27").Value = 60400688
$ws.Range("G27").Value = 75111547
$ws.Range("H27").Value = 78724508

$ws.Range("D29").Value = 7097904
$ws.Range("E29").Value = 6353021
$ws.Range("F29").Value = 7858641
$ws.Range("G29").Value = 19344758
$ws.Range("H29").Value = 2160672

$ws.Range("D30").Value = "-"
$ws.Range("E30").Value = "-"
$ws.Range("F30").Value = "-"
$ws.Range("G30").Value = "-"
$ws.Range("H30").Value = "-"

$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("F31").Value = 0
$ws.Range("G31").Value = 0
$ws.Range("H31").Value = 959708

$ws.Range("D32").Value = 31493
$ws.Range("E32").Value = 29204
$ws.Range("F32").Value = 29107
$ws.Range("G32").Value = 3184785
$ws.Range("H32").Value = 10084822

$ws.Range("D33").Value = 98572
$ws.Range("E33").Value = 105536
$ws.Range("F33").Value = 119846
$ws.Range("G33").Value = 418884
$ws.Range("H33").Value = 458385

$ws.Range("D34").Value = 526346
$ws.Range("E34").Value = 821108
$ws.Range("F34").Value = 1424821
$ws.Range("G34").Value = 2666667
$ws.Range("H34").Value = 0

$ws.Range("D35").Value = 1559590
$ws.Range("E35").Value = 896925
$ws.Range("F35").Value = 928286
$ws.Range("G35").Value = 1042900
$ws.Range("H35").Value = 6508429

$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 0

$ws.Range("D37").Value = 9313905
$ws.Range("E37").Value = 8205794
$ws.Range("F37").Value = 10360701
$ws.Range("G37").Value = 26657994
$ws.Range("H37").Value = 20172016

$ws.Range("D38").Value = 425357
$ws.Range("E38").Value = 0
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 0
$ws.Range("H38").Value = 0

$ws.Range("D39").Value = "-"
$ws.Range("E39").Value = "-"
$ws.Range("F39").Value = "-"
$ws.Range("G39").Value = "-"
$ws.Range("H39").Value = "-"

$ws.Range("D40").Value = 0
$ws.Range("E40").Value = 0
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 0
$ws.Range("H40").Value = 0

$ws.Range("D41").Value = 626073
$ws.Range("E41").Value = 1704293
$ws.Range("F41").Value = 2101859
$ws.Range("G41").Value = 3144095
$ws.Range("H41").Value = 4950335

$ws.Range("D42").Value = 1051430
$ws.Range("E42").Value = 1704293
$ws.Range("F42").Value = 2101859
$ws.Range("G42").Value = 3144095
$ws.Range("H42").Value = 4950335

$ws.Range("D43").Value = 10365335
$ws.Range("E43").Value = 9910087
$ws.Range("F43").Value = 12462560
$ws.Range("G43").Value = 29802089
$ws.Range("H43").Value = 25122351

$ws.Range("D45").Value = 2000000
$ws.Range("E45").Value = 6000000
$ws.Range("F45").Value = 6000000
$ws.Range("G45").Value = 6000000
$ws.Range("H45").Value = 6000000

$ws.Range("D46").Value = 0
$ws.Range("E46").Value = 0
$ws.Range("F46").Value = 0
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 0

$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0

$ws.Range("D48").Value = 0
$ws.Range("E48").Value = 0
$ws.Range("F48").Value = -115975
$ws.Range("G48").Value = -70240
$ws.Range("H48").Value = -334248

$ws.Range("D49").Value = 0
$ws.Range("E49").Value = 0
$ws.Range("F49").Value = 4656
$ws.Range("G49").Value = 35925
$ws.Range("H49").Value = 118611

$ws.Range("D50").Value = 200000
$ws.Range("E50").Value = 600000
$ws.Range("F50").Value = 600000
$ws.Range("G50").Value = 600000
$ws.Range("H50").Value = 600000

$ws.Range("D51").Value = 50000
$ws.Range("E51").Value = 50000
$ws.Range("F51").Value = 50000
$ws.Range("G51").Value = 50000
$ws.Range("H51").Value = 50000

$ws.Range("D52").Value = "-"
$ws.Range("E52").Value = "-"
$ws.Range("F52").Value = "-"
$ws.Range("G52").Value = "-"
$ws.Range("H52").Value = "-"

$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("F53").Value = 0
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 0

$ws.Range("D54").Value = "-"
$ws.Range("E54").Value = "-"
$ws.Range("F54").Value = "-"
$ws.Range("G54").Value = "-"
$ws.Range("H54").Value = "-"

$ws.Range("D55").Value = 0
$ws.Range("E55").Value = 0
$ws.Range("F55").Value = 0
$ws.Range("G55").Value = 0
$ws.Range("H55").Value = 0

$ws.Range("D56").Value = 31577951
$ws.Range("E56").Value = 27680704
$ws.Range("F56").Value = 41399447
$ws.Range("G56").Value = 38693773
$ws.Range("H56").Value = 47167794

$ws.Range("D57").Value = 33827951
$ws.Range("E57").Value = 34330704
$ws.Range("F57").Value = 47938128
$ws.Range("G57").Value = 45309458
$ws.Range("H57").Value = 53602157

$ws.Range("D58").Value = 44193286
$ws.Range("E58").Value = 44240791
$ws.Range("F58").Value = 60400688
$ws.Range("G58").Value = 75111547
$ws.Range("H58").Value = 78724508
